$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sort the data rows (A1:B23, header included) in descending order by the
# value column (B), keeping the header row on top.
$sortRange = $ws.Range("A1:B23")
$sortKey = $ws.Range("B1")
$xlDescending = 2
$xlYes = 1
$sortRange.Sort($sortKey, $xlDescending, $null, $null, $xlYes, $null, $xlYes, $xlYes, $false, $null, $null, $xlYes)

# After sorting, the two zero-value language rows ("Russian" and "Uzbek")
# end up at the bottom (rows 22-23). Remove them entirely, shifting the
# remaining cells up so the used range becomes A1:B21.
$xlShiftUp = -4162
$ws.Range("A22:B23").Delete($xlShiftUp)
